$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The "Price" column holds numeric-looking strings (e.g. "230.60", "1.00",
# "2.351.97") that must stay literal text -- Excel would otherwise coerce them
# to real numbers and silently drop trailing zeros / thousand-dot grouping.
# Flip the cell to text format for the assignment, then restore "General" so
# the cell formatting matches the original workbook.
$ws.Range('D2').NumberFormat = '@'
$ws.Range('D2').Value = '37.301.74'
$ws.Range('D2').NumberFormat = 'General'
$ws.Range('E2').Value = '  -1.40%  '
$ws.Range('D3').NumberFormat = '@'
$ws.Range('D3').Value = '2.051.31'
$ws.Range('D3').NumberFormat = 'General'
$ws.Range('E3').Value = '  -1.56%  '
$ws.Range('D4').NumberFormat = '@'
$ws.Range('D4').Value = '0.999'
$ws.Range('D4').NumberFormat = 'General'
$ws.Range('E4').Value = '  -0.17%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '230.60'
$ws.Range('D5').NumberFormat = 'General'
$ws.Range('E5').Value = '  -1.32%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '0.620'
$ws.Range('D6').NumberFormat = 'General'
$ws.Range('E6').Value = '  -0.94%  '
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '56.98'
$ws.Range('D8').NumberFormat = 'General'
$ws.Range('E9').Value = '  -3.08%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '0.0769'
$ws.Range('D10').NumberFormat = 'General'
$ws.Range('E10').Value = '  -2.75%  '
$ws.Range('E11').Value = '  +1.19%  '
$ws.Range('B12').Value = 'WrappedliquidstakedEther2.0'
$ws.Range('C12').Value = 'https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth'
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '2.351.97'
$ws.Range('D12').NumberFormat = 'General'
$ws.Range('E12').Value = '  -1.54%  '
$ws.Range('B13').Value = 'Chainlink'
$ws.Range('C13').Value = 'https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link'
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '14.64'
$ws.Range('D13').NumberFormat = 'General'
$ws.Range('E13').Value = '  -1.25%  '
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '20.57'
$ws.Range('D14').NumberFormat = 'General'
$ws.Range('E14').Value = '  -3.26%  '
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '0.756'
$ws.Range('D15').NumberFormat = 'General'
$ws.Range('E15').Value = '  -2.87%  '
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '5.25'
$ws.Range('D16').NumberFormat = 'General'
$ws.Range('E16').Value = '  -2.02%  '
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '2.053.89'
$ws.Range('D17').NumberFormat = 'General'
$ws.Range('E17').Value = '  -0.56%  '
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '37.246.53'
$ws.Range('D18').NumberFormat = 'General'
$ws.Range('E18').Value = '  -1.38%  '
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '6.04'
$ws.Range('D19').NumberFormat = 'General'
$ws.Range('E19').Value = '  -2.23%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '69.66'
$ws.Range('D20').NumberFormat = 'General'
$ws.Range('E20').Value = '  -2.67%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '0.0₃0821'
$ws.Range('D21').NumberFormat = 'General'
$ws.Range('E21').Value = '  -3.63%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '226.35'
$ws.Range('D22').NumberFormat = 'General'
$ws.Range('E22').Value = '  -0.87%  '
$ws.Range('E23').Value = '  +0.03%  '
$ws.Range('E24').Value = '  +0.55%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '2.33'
$ws.Range('D25').NumberFormat = 'General'
$ws.Range('E25').Value = '  -3.94%  '
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '9.76'
$ws.Range('D26').NumberFormat = 'General'
$ws.Range('E26').Value = '  +6.07%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '169.97'
$ws.Range('D27').NumberFormat = 'General'
$ws.Range('E27').Value = '  -1.10%  '
$ws.Range('E28').Value = '  -6.17%  '
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '19.15'
$ws.Range('D29').NumberFormat = 'General'
$ws.Range('E29').Value = '  -1.82%  '
$ws.Range('E30').Value = '  -5.77%  '
$ws.Range('E31').Value = '  -0.41%  '
$ws.Range('E32').Value = '  -4.75%  '
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '0.0622'
$ws.Range('D33').NumberFormat = 'General'
$ws.Range('E33').Value = '  -1.82%  '
$ws.Range('E34').Value = '  -4.56%  '
$ws.Range('E35').Value = '  -1.49%  '
$ws.Range('E36').Value = '  +0.06%  '
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '3.26'
$ws.Range('D37').NumberFormat = 'General'
$ws.Range('E37').Value = '  -5.49%  '
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '1.00'
$ws.Range('D38').NumberFormat = 'General'
$ws.Range('E38').Value = '  +0.08%  '
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '5.30'
$ws.Range('D39').NumberFormat = 'General'
$ws.Range('E39').Value = '  -2.33%  '
$ws.Range('E40').Value = '  +3.31%  '
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '98.25'
$ws.Range('D41').NumberFormat = 'General'
$ws.Range('E41').Value = '  -1.15%  '
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '1.481.12'
$ws.Range('D42').NumberFormat = 'General'
$ws.Range('E42').Value = '  +2.42%  '
$ws.Range('E43').Value = '  +0.46%  '
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '0.0950'
$ws.Range('D44').NumberFormat = 'General'
$ws.Range('E44').Value = '  -3.64%  '
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '1.17'
$ws.Range('D45').NumberFormat = 'General'
$ws.Range('E45').Value = '  +1.04%  '
$ws.Range('E46').Value = '  -1.79%  '
$ws.Range('B47').Value = 'ARBITRUM'
$ws.Range('C47').Value = 'https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb'
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '1.03'
$ws.Range('D47').NumberFormat = 'General'
$ws.Range('E47').Value = '  -3.50%  '
$ws.Range('B48').Value = 'FTXToken'
$ws.Range('C48').Value = 'https://coinranking.com/coin/NfeOYfNcl+ftxtoken-ftt'
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '3.97'
$ws.Range('D48').NumberFormat = 'General'
$ws.Range('E48').Value = '  -5.34%  '
$ws.Range('E50').Value = '  -2.33%  '
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '2.237.70'
$ws.Range('D51').NumberFormat = 'General'
$ws.Range('E51').Value = '  -1.58%  '
